# Insert a new data row at row 243 (shifting the existing rows 243:340 down
# to 244:341) and populate it with a new observation for
# Hortaliza / Femacal de La Calera - Berenjena.
#
# The new row mirrors the market stats of the (now shifted) row that used to
# sit at 269 (Volumen=105, Precio min/max/prom=9000/9500/9262, Precio $/Kg=154,
# Origen="Región de Arica y Parinacota"), but is stamped with a newer date
# (2022-08-22, Excel serial 44795).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 243..340 down to 244..341, leaving row 243 blank.
$ws.Rows.Item(243).EntireRow.Insert()

$ws.Range("A243").Value = 3
$ws.Range("B243").Value = "Femacal de La Calera"
$ws.Range("C243").Value = "Coquimbo"
$ws.Range("D243").Value = 44795
$ws.Range("E243").Value = 5
$ws.Range("F243").Value = 100112001
$ws.Range("G243").Value = "Berenjena"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 105
$ws.Range("K243").Value = 9000
$ws.Range("L243").Value = 9500
$ws.Range("M243").Value = 9262
$ws.Range("N243").Value = "`$/caja 60 unidades"
$ws.Range("O243").Value = "Región de Arica y Parinacota"
$ws.Range("P243").Value = 154
$ws.Range("Q243").Value = 60
$ws.Range("R243").Value = "Hortaliza"
